$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1481
$ws1.Range("F3").Value = 1448
$ws1.Range("F5").Value = 228
$ws1.Range("F6").Value = 721
$ws1.Range("F7").Value = 39
$ws1.Range("F8").Value = 645
$ws1.Range("F11").Value = 1385
$ws1.Range("F12").Value = 34088
$ws1.Range("F13").Value = 7141
$ws1.Range("F14").Value = 116
$ws1.Range("F15").Value = 373
$ws1.Range("F16").Value = 584
$ws1.Range("F17").Value = 449
$ws1.Range("F19").Value = 110
$ws1.Range("F20").Value = 181
$ws1.Range("F21").Value = 51
$ws1.Range("F23").Value = 109
$ws1.Range("F24").Value = 812
$ws1.Range("F25").Value = 16
$ws1.Range("F26").Value = 322
$ws1.Range("F27").Value = 392
$ws1.Range("F28").Value = 445
$ws1.Range("F30").Value = 215
$ws1.Range("F31").Value = 54
$ws1.Range("F32").Value = 744
$ws1.Range("F35").Value = 755
$ws1.Range("F36").Value = 115
$ws1.Range("F38").Value = 801
$ws1.Range("F39").Value = 291
$ws1.Range("F41").Value = 26

# --- Sheet 2: 演出 (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 1213
$ws2.Range("F5").Value = 163
$ws2.Range("F6").Value = 292
$ws2.Range("F7").Value = 4329
$ws2.Range("F9").Value = 241
$ws2.Range("F13").Value = 40
$ws2.Range("F14").Value = 6
$ws2.Range("F19").Value = 4302

# --- Sheet 3: 本地生活 (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 1468
$ws3.Range("F3").Value = 360

# --- Sheet 4: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")

# --- Sheet 4 rows 13-15: replace event listings (shifted by new event insertion upstream) ---
$ws4.Range("B13").NumberFormat = "@"
$ws4.Range("B13").Value = "2024-07-20"
$ws4.Range("C13").Value = "广州·冰兔2024线下live「过去和未来」"
$ws4.Range("D13").Value = "恩宁路265号三层四层自编01 MAO Livehouse广州（永庆坊店）"
$ws4.Range("E13").Value = "2024.07.20 20:00-07.20 22:00"
$ws4.Range("F13").Value = 163
$ws4.Range("G13").Value = 198
$ws4.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=87546"
$ws4.Range("I13").Value = "//i2.hdslb.com/bfs/openplatform/202406/2X09PE1a1718611339266.jpeg"

$ws4.Range("C14").Value = "广州·跨越二次元ACG神级动漫世界巡回演唱会"
$ws4.Range("D14").Value = "广州市荔湾区十甫路125号(上下九步行街内)2层 广州平安大戏院"
$ws4.Range("E14").Value = "2024.07.20 19:30-07.20 21:10"
$ws4.Range("F14").Value = 292
$ws4.Range("G14").Value = 280
$ws4.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=85353"
$ws4.Range("I14").Value = "//i1.hdslb.com/bfs/openplatform/202405/4gACWbPh1715223804704.jpeg"

$ws4.Range("C15").Value = "广州·音波狂潮II 萤光宇宙 音游嘉年华"
$ws4.Range("D15").Value = "新港东路磨碟沙大街118号自编8栋 啤厂媒棚"
$ws4.Range("E15").Value = "2024.07.20 13:30-07.21 23:30"
$ws4.Range("F15").Value = 4329
$ws4.Range("G15").Value = "已售罄"
$ws4.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=86632"
$ws4.Range("I15").Value = "//i0.hdslb.com/bfs/openplatform/202405/GcK1JV3B1717123497026.jpeg"

$ws4.Range("F2").Value = 1468
$ws4.Range("F3").Value = 360
$ws4.Range("F4").Value = 1213
$ws4.Range("F5").Value = 1481
$ws4.Range("F7").Value = 1448
$ws4.Range("F8").Value = 228
$ws4.Range("F9").Value = 721
$ws4.Range("F10").Value = 39
$ws4.Range("F11").Value = 645
$ws4.Range("F16").Value = 241
$ws4.Range("F17").Value = 241
$ws4.Range("F20").Value = 7142
$ws4.Range("F21").Value = 373
$ws4.Range("F23").Value = 584
$ws4.Range("F24").Value = 449
$ws4.Range("F25").Value = 40
$ws4.Range("F26").Value = 110
$ws4.Range("F27").Value = 182
$ws4.Range("F28").Value = 6
$ws4.Range("F29").Value = 51
$ws4.Range("F32").Value = 109
$ws4.Range("F33").Value = 812
$ws4.Range("F34").Value = 16
$ws4.Range("F35").Value = 322
$ws4.Range("F36").Value = 392
$ws4.Range("F37").Value = 445
$ws4.Range("F39").Value = 215
$ws4.Range("F40").Value = 54
$ws4.Range("F41").Value = 744
$ws4.Range("F43").Value = 293
$ws4.Range("F45").Value = 801
$ws4.Range("F46").Value = 291
$ws4.Range("F49").Value = 26
